# Updates cryptos list (coinranking.com snapshot) with refreshed prices/volumes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.675.92'
$ws.Range("E2").Value = '  +1.23%  '

# Row 3
$ws.Range("D3").Value = '1.692.52'
$ws.Range("E3").Value = '  +0.34%  '

# Row 4
$ws.Range("E4").Value = '  -0.59%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.36'
$ws.Range("E5").Value = '  +1.06%  '

# Row 6
$ws.Range("E6").Value = '  -0.37%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3944'
$ws.Range("E7").Value = '  +0.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4058'
$ws.Range("E8").Value = '  +0.95%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.491'
$ws.Range("E9").Value = '  -0.33%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.003'
$ws.Range("E10").Value = '  -0.63%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.20'
$ws.Range("E11").Value = '  -2.72%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08870'
$ws.Range("E12").Value = '  +1.65%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.253'
$ws.Range("E13").Value = '  +0.03%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.64'
$ws.Range("E14").Value = '  +3.08%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.055'
$ws.Range("E15").Value = '  +7.85%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001319'
$ws.Range("E16").Value = '  +0.53%  '

# Row 17
$ws.Range("D17").Value = '1.691.88'
$ws.Range("E17").Value = '  -0.41%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.69'
$ws.Range("E18").Value = '  -0.33%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07022'
$ws.Range("E19").Value = '  -1.09%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.57'
$ws.Range("E20").Value = '  +1.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.001'
$ws.Range("E21").Value = '  +4.86%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  -0.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.38'
$ws.Range("E23").Value = '  +2.60%  '

# Row 24
$ws.Range("D24").Value = '24.665.20'
$ws.Range("E24").Value = '  +1.20%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.260'
$ws.Range("E25").Value = '  +10.66%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.357'
$ws.Range("E26").Value = '  +1.52%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.77'
$ws.Range("E27").Value = '  +2.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.44'
$ws.Range("E28").Value = '  +2.37%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '136.12'
$ws.Range("E29").Value = '  +2.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.183'
$ws.Range("E30").Value = '  +0.73%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.619'
$ws.Range("E31").Value = '  +3.66%  '

# Row 32
$ws.Range("D32").Value = '1.879.20'
$ws.Range("E32").Value = '  -0.40%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08600'
$ws.Range("E33").Value = '  +0.43%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.058'
$ws.Range("E34").Value = '  -1.49%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.116'
$ws.Range("E35").Value = '  -2.22%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.28'
$ws.Range("E36").Value = '  +2.36%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2734'
$ws.Range("E37").Value = '  +1.73%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.886'
$ws.Range("E38").Value = '  -4.15%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.50'
$ws.Range("E39").Value = '  -0.59%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09218'
$ws.Range("E40").Value = '  +3.16%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02722'
$ws.Range("E41").Value = '  -0.82%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.467'
$ws.Range("E42").Value = '  +0.29%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7658'
$ws.Range("E43").Value = '  +1.20%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.13'
$ws.Range("E44").Value = '  +4.60%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.596'
$ws.Range("E45").Value = '  +7.20%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7150'
$ws.Range("E46").Value = '  +0.97%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.217'
$ws.Range("E47").Value = '  +1.40%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  -0.25%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.13'
$ws.Range("E49").Value = '  +0.97%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.321'
$ws.Range("E50").Value = '  +3.21%  '

# Row 51 - coin replaced (Cronos -> Aave) with new price/volume
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '90.60'
$ws.Range("E51").Value = '  +5.73%  '

